$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Ingredientes": replace the "pão sírio" ingredient list with the
# "bolo de banana" ingredient list.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Ingredientes")

# Drop the old rows 10-15 (tomato/onion/lime/meat/"sal a gosto" rows) first,
# leaving rows 2-9 to be overwritten in place.
$ws1.Rows("10:15").Delete()

$ingredientes = @(
    @(3,   "bananas maduras amassadas"),
    @(2,   "ovos"),
    @(0.5, "xícara de açúcar"),
    @(0.5, "xícara de óleo vegetal"),
    @(1,   "xícara de aveia em flocos"),
    @(1,   "xícara de farinha de trigo"),
    @(1,   "colher de sopa de canela em pó"),
    @(1,   "colher de sopa de fermento em pó")
)

for ($i = 0; $i -lt $ingredientes.Count; $i++) {
    $row = 2 + $i
    $ws1.Range("A$row").Value = $ingredientes[$i][0]
    $ws1.Range("B$row").Value = $ingredientes[$i][1]
}

$ws1.Range("D2:D9").Formula = '="<li><span class=""qtd"" data-qtd="""&A2&""">"&A2&"</span> "&B2&"</li>"'

# ---------------------------------------------------------------------------
# Sheet "Preparo": replace the "pão sírio" preparation steps with the
# "bolo de banana" steps.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Preparo")

$ws2.Rows("8:10").Delete()

$passos = @(
    "Pré-aqueça o forno a 180°C e unte uma forma com margarina e farinha de trigo.",
    "Em uma tigela grande, misture as bananas amassadas, os ovos, o açúcar e o óleo até obter uma mistura homogênea.",
    "Adicione a aveia, a farinha de trigo e a canela em pó e misture bem.",
    "Adicione o fermento em pó e mexa suavemente até incorporar à massa.",
    "Despeje a massa na forma untada e leve ao forno por cerca de 30 a 40 minutos, ou até que o bolo esteja dourado e assado.",
    "Retire do forno e deixe esfriar antes de servir."
)

for ($i = 0; $i -lt $passos.Count; $i++) {
    $row = 2 + $i
    $ws2.Range("A$row").Value = $passos[$i]
}

# Matches the author's original edit pattern: rows 2 and 3 end up with their
# own (unshared) formula, while rows 4-7 form a new shared-formula group.
$ws2.Range("B2").Formula = '="<li>"&A2&"</li>"'
$ws2.Range("B3").Formula = '="<li>"&A3&"</li>"'
$ws2.Range("B4:B7").Formula = '="<li>"&A4&"</li>"'
